# Apply "Fruta / hortaliza, semanal" weekly refresh to rows 4-11 and 13-14.
# Row data is rotated between the date/quality/price columns (D, K, L, M, N, O, P, Q, R, S).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 <- previous row 13 data
$ws.Range("D4").Value = 44243
$ws.Range("K4").Value = "Black Amber"
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 300
$ws.Range("N4").Value = 14000
$ws.Range("O4").Value = 15000
$ws.Range("P4").Value = 14500
$ws.Range("Q4").Value = "$/caja 18 kilos granel"
$ws.Range("R4").Value = "Región de O'Higgins"
$ws.Range("S4").Value = 806

# Row 5 <- previous row 10 data
$ws.Range("D5").Value = 44174
$ws.Range("K5").Value = "Angeleno"
$ws.Range("L5").Value = "Primera"
$ws.Range("M5").Value = 270
$ws.Range("N5").Value = 20000
$ws.Range("O5").Value = 21000
$ws.Range("P5").Value = 20500
$ws.Range("Q5").Value = "$/caja 18 kilos granel"
$ws.Range("R5").Value = "Región de O'Higgins"
$ws.Range("S5").Value = 1139

# Row 6 <- previous row 9 data
$ws.Range("D6").Value = 44245
$ws.Range("K6").Value = "Black Amber"
$ws.Range("L6").Value = "Primera"
$ws.Range("M6").Value = 250
$ws.Range("N6").Value = 14000
$ws.Range("O6").Value = 15000
$ws.Range("P6").Value = 14500
$ws.Range("Q6").Value = "$/bandeja 18 kilos granel"
$ws.Range("R6").Value = "Región de O'Higgins"
$ws.Range("S6").Value = 806

# Row 7 <- previous row 5 data
$ws.Range("D7").Value = 44238
$ws.Range("K7").Value = "Black Amber"
$ws.Range("L7").Value = "Segunda"
$ws.Range("M7").Value = 300
$ws.Range("N7").Value = 14000
$ws.Range("O7").Value = 15000
$ws.Range("P7").Value = 14500
$ws.Range("Q7").Value = "$/bandeja 18 kilos granel"
$ws.Range("R7").Value = "Región de O'Higgins"
$ws.Range("S7").Value = 806

# Row 8 <- previous row 6 data
$ws.Range("D8").Value = 44238
$ws.Range("K8").Value = "Fortuna"
$ws.Range("L8").Value = "Segunda"
$ws.Range("M8").Value = 300
$ws.Range("N8").Value = 14000
$ws.Range("O8").Value = 15000
$ws.Range("P8").Value = 14500
$ws.Range("Q8").Value = "$/bandeja 18 kilos granel"
$ws.Range("R8").Value = "Región de O'Higgins"
$ws.Range("S8").Value = 806

# Row 9 <- previous row 11 data
$ws.Range("D9").Value = 44169
$ws.Range("K9").Value = "Angeleno"
$ws.Range("L9").Value = "Tercera"
$ws.Range("M9").Value = 250
$ws.Range("N9").Value = 24000
$ws.Range("O9").Value = 25000
$ws.Range("P9").Value = 24500
$ws.Range("Q9").Value = "$/bandeja 18 kilos granel"
$ws.Range("R9").Value = "Región de O'Higgins"
$ws.Range("S9").Value = 1361

# Row 10 <- previous row 14 data
$ws.Range("D10").Value = 44314
$ws.Range("K10").Value = "Angeleno"
$ws.Range("L10").Value = "Segunda"
$ws.Range("M10").Value = 250
$ws.Range("N10").Value = 14000
$ws.Range("O10").Value = 15000
$ws.Range("P10").Value = 14500
$ws.Range("Q10").Value = "$/bandeja 18 kilos granel"
$ws.Range("R10").Value = "Región de O'Higgins"
$ws.Range("S10").Value = 806

# Row 11 <- previous row 8 data
$ws.Range("D11").Value = 44278
$ws.Range("K11").Value = "Angeleno"
$ws.Range("L11").Value = "Primera"
$ws.Range("M11").Value = 300
$ws.Range("N11").Value = 15000
$ws.Range("O11").Value = 16000
$ws.Range("P11").Value = 15500
$ws.Range("Q11").Value = "$/caja 18 kilos granel"
$ws.Range("R11").Value = "Región de O'Higgins"
$ws.Range("S11").Value = 861

# Row 13 <- previous row 7 data
$ws.Range("D13").Value = 44175
$ws.Range("K13").Value = "Angeleno"
$ws.Range("L13").Value = "Primera"
$ws.Range("M13").Value = 200
$ws.Range("N13").Value = 21000
$ws.Range("O13").Value = 22000
$ws.Range("P13").Value = 21500
$ws.Range("Q13").Value = "$/bandeja 18 kilos granel"
$ws.Range("R13").Value = "Región de O'Higgins"
$ws.Range("S13").Value = 1194

# Row 14 <- previous row 4 data
$ws.Range("D14").Value = 44217
$ws.Range("K14").Value = "Black Amber"
$ws.Range("L14").Value = "Segunda"
$ws.Range("M14").Value = 300
$ws.Range("N14").Value = 16000
$ws.Range("O14").Value = 17000
$ws.Range("P14").Value = 16500
$ws.Range("Q14").Value = "$/bandeja 18 kilos granel"
$ws.Range("R14").Value = "Región Metropolitana"
$ws.Range("S14").Value = 917

